$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 4.077281088115026
$ws.Cells.Item(2, 5).Value = 3.013456857613561
$ws.Cells.Item(3, 3).Value = 4.731933890736406
$ws.Cells.Item(3, 5).Value = 4.482374336309247
$ws.Cells.Item(4, 3).Value = 7.427427238257622
$ws.Cells.Item(4, 5).Value = 6.028771162940338
$ws.Cells.Item(5, 3).Value = 3.697244903694119
$ws.Cells.Item(5, 5).Value = 5.35733476656457
$ws.Cells.Item(6, 3).Value = -0.9010546343133807
$ws.Cells.Item(6, 5).Value = 1.483308540745609
$ws.Cells.Item(7, 3).Value = 1.278852728916302
$ws.Cells.Item(7, 5).Value = 1.211545622441634
$ws.Cells.Item(8, 3).Value = 2.40113223806655
$ws.Cells.Item(8, 5).Value = 1.221134982408678
$ws.Cells.Item(9, 3).Value = 0.4170416928886977
$ws.Cells.Item(9, 5).Value = 1.535040327807513
$ws.Cells.Item(10, 3).Value = 2.025199397970145
$ws.Cells.Item(10, 5).Value = 1.766027057877517
$ws.Cells.Item(11, 3).Value = 2.535130037318867
$ws.Cells.Item(11, 5).Value = 2.27519030414034
$ws.Cells.Item(12, 3).Value = 1.085017960020163
$ws.Cells.Item(12, 5).Value = 1.594198396297974
$ws.Cells.Item(13, 3).Value = 2.388449447315399
$ws.Cells.Item(13, 5).Value = 1.896063224966515
$ws.Cells.Item(14, 3).Value = 1.940295589655605
$ws.Cells.Item(14, 5).Value = 2.152498180268503
$ws.Cells.Item(15, 3).Value = 0.1294783794713039
$ws.Cells.Item(15, 5).Value = 1.357330623126884
$ws.Cells.Item(16, 3).Value = 0.2239492534813481
$ws.Cells.Item(16, 5).Value = 0.1832633470655098
$ws.Cells.Item(17, 3).Value = 0.669497318959178
$ws.Cells.Item(17, 5).Value = 0.4274836977099516
$ws.Cells.Item(18, 3).Value = 0.9511721486624936
$ws.Cells.Item(18, 5).Value = 0.886066450790457
$ws.Cells.Item(19, 3).Value = 1.350833417525776
$ws.Cells.Item(19, 5).Value = 1.249180524815863
$ws.Cells.Item(20, 3).Value = 3.305550968939119
$ws.Cells.Item(20, 5).Value = 2.444559947892744
$ws.Cells.Item(21, 3).Value = 2.741128804567849
$ws.Cells.Item(21, 5).Value = 3.207064487734335
$ws.Cells.Item(22, 3).Value = -5.478868953971427
$ws.Cells.Item(22, 5).Value = -2.015486574969738
$ws.Cells.Item(23, 3).Value = -0.1094048593225039
$ws.Cells.Item(23, 5).Value = -1.306092631642397
$ws.Cells.Item(24, 3).Value = 3.381937564063731
$ws.Cells.Item(24, 5).Value = 0.8627271536207459
$ws.Cells.Item(25, 3).Value = 1.637009187238481
$ws.Cells.Item(25, 5).Value = 2.774647569643585
$ws.Cells.Item(26, 3).Value = -0.03183845066089264
$ws.Cells.Item(26, 5).Value = 1.102201969172678
$ws.Cells.Item(27, 3).Value = 1.810762846774527
$ws.Cells.Item(27, 5).Value = 1.207964354105195
$ws.Cells.Item(28, 3).Value = 1.292027346513414
$ws.Cells.Item(28, 5).Value = 0.9897401519578963
$ws.Cells.Item(29, 3).Value = 1.566699735746391
$ws.Cells.Item(29, 5).Value = 1.595944879287448
$ws.Cells.Item(30, 3).Value = 1.752318341645176
$ws.Cells.Item(30, 5).Value = 1.832607040765044
$ws.Cells.Item(31, 3).Value = 2.349173111882341
$ws.Cells.Item(31, 5).Value = 2.133017022807637
$ws.Cells.Item(32, 3).Value = 0.8252516928923814
$ws.Cells.Item(32, 5).Value = 1.704160690624001
$ws.Cells.Item(33, 3).Value = -1.611564731980897
$ws.Cells.Item(33, 5).Value = -0.741724804865096
$ws.Cells.Item(34, 3).Value = -0.5453046728852495
$ws.Cells.Item(34, 5).Value = -1.767889269204159
$ws.Cells.Item(35, 3).Value = 1.726013280798222
$ws.Cells.Item(35, 5).Value = 0.05925426804285205
$ws.Cells.Item(36, 3).Value = -0.06641493770841445
$ws.Cells.Item(36, 5).Value = 0.7991555368092929
$ws.Cells.Item(37, 3).Value = -0.1397017661237232
$ws.Cells.Item(37, 5).Value = 0.385156833908451
$ws.Cells.Item(38, 3).Value = 0.06464796496492564
$ws.Cells.Item(38, 5).Value = -0.06351196001971315
